$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-15 21:48:57"
$ws.Range("E3").Value = "2026-02-15 21:49:00"
$ws.Range("G3").Value = "200 cm"
$ws.Range("H3").Value = "97%"
$ws.Range("I3").Value = "2.6 mm"
$ws.Range("O3").Value = "-4.9 °C"
$ws.Range("E4").Value = "2026-02-15 21:49:03"
$ws.Range("O4").Value = "7.3 °C"
$ws.Range("E5").Value = "2026-02-15 21:49:05"
$ws.Range("I5").Value = "8.0 mm"
$ws.Range("O5").Value = "-4.3 °C"
$ws.Range("E6").Value = "2026-02-15 21:49:08"
$ws.Range("J6").Value = "1015.6 hPa"
$ws.Range("E7").Value = "2026-02-15 21:49:11"
$ws.Range("E8").Value = "2026-02-15 21:49:13"
$ws.Range("E9").Value = "2026-02-15 21:49:16"
$ws.Range("H9").Value = "53%"
$ws.Range("E10").Value = "2026-02-15 21:49:19"
$ws.Range("H10").Value = "72%"
$ws.Range("E11").Value = "2026-02-15 21:49:21"
$ws.Range("H11").Value = "46%"
$ws.Range("O11").Value = "7.0 °C"
$ws.Range("E12").Value = "2026-02-15 21:49:24"
$ws.Range("H12").Value = "59%"
$ws.Range("O12").Value = "10.7 °C"
$ws.Range("E13").Value = "2026-02-15 21:49:26"
$ws.Range("H13").Value = "39%"
$ws.Range("J13").Value = "1015.5 hPa"
$ws.Range("O13").Value = "6.3 °C"
$ws.Range("E14").Value = "2026-02-15 21:49:29"
$ws.Range("K14").Value = "11.5 MJ/m2"
$ws.Range("E15").Value = "2026-02-15 21:49:32"
$ws.Range("H15").Value = "53%"
$ws.Range("E16").Value = "2026-02-15 21:49:34"
$ws.Range("H16").Value = "63%"
$ws.Range("E17").Value = "2026-02-15 21:49:37"
$ws.Range("E18").Value = "2026-02-15 21:49:40"
$ws.Range("E19").Value = "2026-02-15 21:49:42"
$ws.Range("E20").Value = "2026-02-15 21:49:45"
$ws.Range("L20").Value = "80.6 km/h - 332º 21:27 TU"
$ws.Range("E21").Value = "2026-02-15 21:49:48"
$ws.Range("H21").Value = "40%"
$ws.Range("J21").Value = "1015.1 hPa"
$ws.Range("E22").Value = "2026-02-15 21:49:51"
$ws.Range("E23").Value = "2026-02-15 21:49:53"
$ws.Range("H23").Value = "67%"
$ws.Range("I23").Value = "4.9 mm"
$ws.Range("O23").Value = "-3.4 °C"
$ws.Range("E24").Value = "2026-02-15 21:49:56"
$ws.Range("O24").Value = "9.0 °C"
$ws.Range("E25").Value = "2026-02-15 21:49:59"
$ws.Range("O25").Value = "-1.3 °C"
$ws.Range("E26").Value = "2026-02-15 21:50:02"
$ws.Range("E27").Value = "2026-02-15 21:50:04"
$ws.Range("O27").Value = "0.1 °C"
$ws.Range("E28").Value = "2026-02-15 21:50:07"
$ws.Range("O28").Value = "6.5 °C"
$ws.Range("E29").Value = "2026-02-15 21:50:09"
$ws.Range("H29").Value = "60%"
$ws.Range("E30").Value = "2026-02-15 21:50:12"
$ws.Range("O30").Value = "9.7 °C"
$ws.Range("E31").Value = "2026-02-15 21:50:15"
$ws.Range("E32").Value = "2026-02-15 21:50:18"
$ws.Range("E33").Value = "2026-02-15 21:50:20"
$ws.Range("H33").Value = "43%"
$ws.Range("J33").Value = "1015.1 hPa"
$ws.Range("L33").Value = "33.5 km/h - 269º 21:06 TU"
$ws.Range("O33").Value = "6.0 °C"
$ws.Range("E34").Value = "2026-02-15 21:50:23"
$ws.Range("H34").Value = "52%"
$ws.Range("O34").Value = "1.3 °C"
$ws.Range("E35").Value = "2026-02-15 21:50:26"
$ws.Range("E36").Value = "2026-02-15 21:50:28"
$ws.Range("H36").Value = "50%"
$ws.Range("O36").Value = "11.3 °C"
$ws.Range("E37").Value = "2026-02-15 21:50:31"
$ws.Range("H37").Value = "56%"
$ws.Range("J37").Value = "1016.5 hPa"
$ws.Range("O37").Value = "5.9 °C"
$ws.Range("E38").Value = "2026-02-15 21:50:34"
$ws.Range("E39").Value = "2026-02-15 21:50:37"
$ws.Range("H39").Value = "60%"
$ws.Range("O39").Value = "-2.6 °C"
$ws.Range("E40").Value = "2026-02-15 21:50:39"
$ws.Range("H40").Value = "41%"
$ws.Range("O40").Value = "8.6 °C"
$ws.Range("E41").Value = "2026-02-15 21:50:42"
$ws.Range("E42").Value = "2026-02-15 21:50:45"
$ws.Range("H42").Value = "59%"
$ws.Range("O42").Value = "10.5 °C"
$ws.Range("E43").Value = "2026-02-15 21:50:47"
$ws.Range("O43").Value = "6.5 °C"
$ws.Range("E44").Value = "2026-02-15 21:50:50"
$ws.Range("I44").Value = "5.5 mm"
$ws.Range("O44").Value = "-3.8 °C"
$ws.Range("E45").Value = "2026-02-15 21:50:53"
$ws.Range("I45").Value = "3.8 mm"
$ws.Range("O45").Value = "1.2 °C"
$ws.Range("E46").Value = "2026-02-15 21:50:56"
